$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix row 3: cedula becomes a real number, and status becomes "Activo" ---
$ws.Range("C3").Value = 123456789
$ws.Range("F3").Value = "Activo"

# --- 2. Add row 4 (new form response), same student as row 3 but a second course ---
$ws.Range("A4").Value = "Luis"
$ws.Range("B4").Value = "Zambrano "
$ws.Range("C4").Value = 123456789
$ws.Range("D4").Value = "jdgaleas1@espe.edu.ec"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:jdgaleas1@espe.edu.ec")
$ws.Range("E4").Value = "Parvularia"
$ws.Range("F4").Value = "Activo"
$ws.Range("G4").Value = "instructorDefaul@gmail.com"
$ws.Hyperlinks.Add($ws.Range("G4"), "mailto:instructorDefaul@gmail.com")

# --- 3. Center-align the response block and apply a zero-padded ID format to column C ---
$ws.Range("A2:G7").HorizontalAlignment = -4108
$ws.Range("A8:G17").HorizontalAlignment = -4108
$ws.Range("C2:C17").NumberFormat = "0000000000"

# --- 4. Stray formatted cell mirrored onto column H, like the existing one on column J ---
$ws.Range("H80").Font.Underline = 2

Write-Host "done"
